$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8681.857
$ws.Range("I62").Value = 10100.5
$ws.Range("J62").Value = 7808.846
$ws.Range("K62").Value = 10100.5
$ws.Range("L62").Value = 7808.846
$ws.Range("M62").Value = -9476.5
$ws.Range("N62").Value = -9056.846
$ws.Range("H65").Value = 8681.857
$ws.Range("I65").Value = 10100.5
$ws.Range("J65").Value = 7808.846
$ws.Range("K65").Value = 50502.5
$ws.Range("L65").Value = 39044.23
$ws.Range("M65").Value = -47382.5
$ws.Range("N65").Value = -45284.23
$ws.Range("H98").Value = 3563.4
$ws.Range("I98").Value = 4276.6665
$ws.Range("K98").Value = 4276.6665
$ws.Range("M98").Value = -2778.6665
$ws.Range("H106").Value = 45455510
$ws.Range("I106").Value = 55556096
$ws.Range("K106").Value = 55556096
$ws.Range("M106").Value = -55555465
$ws.Range("H111").Value = 3834011.5
$ws.Range("I111").Value = 5292695
$ws.Range("K111").Value = 15878085
$ws.Range("M111").Value = -15875018
$ws.Range("H122").Value = 3563.4
$ws.Range("I122").Value = 4276.6665
$ws.Range("K122").Value = 12829.9995
$ws.Range("M122").Value = -10379.9995
$ws.Range("H137").Value = 87930.71000000001
$ws.Range("I137").Value = 120503.734
$ws.Range("J137").Value = 6498.1665
$ws.Range("K137").Value = 361511.202
$ws.Range("L137").Value = 19494.4995
$ws.Range("M137").Value = -358961.202
$ws.Range("N137").Value = -24594.4995
$ws.Range("H138").Value = 5070.8
$ws.Range("J138").Value = 5698.387
$ws.Range("L138").Value = 17095.161
$ws.Range("N138").Value = -27375.161
$ws.Range("H141").Value = 15128.125
$ws.Range("I141").Value = 22008.4
$ws.Range("K141").Value = 66025.20000000001
$ws.Range("M141").Value = -60845.20000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1664771
$ws.Range("I2").Value = 2176435.2
$ws.Range("J2").Value = 1862.25
$ws.Range("K2").Value = 2176435.2
$ws.Range("L2").Value = 1862.25
$ws.Range("M2").Value = -2176322.2
$ws.Range("N2").Value = -2088.25
$ws.Range("H45").Value = 8931943
$ws.Range("I45").Value = 14286917
$ws.Range("J45").Value = 6987.5
$ws.Range("K45").Value = 14286917
$ws.Range("L45").Value = 6987.5
$ws.Range("M45").Value = -14286540
$ws.Range("N45").Value = -7741.5
$ws.Range("H102").Value = 3970648.5
$ws.Range("I102").Value = 4169153.5
$ws.Range("J102").Value = 549
$ws.Range("K102").Value = 4169153.5
$ws.Range("L102").Value = 549
$ws.Range("M102").Value = -4167531.5
$ws.Range("N102").Value = -3793
$ws.Range("H116").Value = 1664771
$ws.Range("I116").Value = 2176435.2
$ws.Range("J116").Value = 1862.25
$ws.Range("K116").Value = 2176435.2
$ws.Range("L116").Value = 1862.25
$ws.Range("M116").Value = -2174141.2
$ws.Range("N116").Value = -6450.25
$ws.Range("H122").Value = 14085147
$ws.Range("I122").Value = 21145944
$ws.Range("J122").Value = 2611350.5
$ws.Range("K122").Value = 63437832
$ws.Range("L122").Value = 7834051.5
$ws.Range("M122").Value = -63435382
$ws.Range("N122").Value = -7838951.5
$ws.Range("H132").Value = 31030.912
$ws.Range("I132").Value = 1781.5518
$ws.Range("K132").Value = 5344.6554
$ws.Range("M132").Value = -2814.6554

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1664771
$ws.Range("I3").Value = 2176435.2
$ws.Range("J3").Value = 1862.25
$ws.Range("K3").Value = 2176435.2
$ws.Range("L3").Value = 1862.25
$ws.Range("M3").Value = -2176321.2
$ws.Range("N3").Value = -2090.25
$ws.Range("H105").Value = 6946981.5
$ws.Range("I105").Value = 8931260
$ws.Range("J105").Value = 2008.5
$ws.Range("K105").Value = 8931260
$ws.Range("L105").Value = 2008.5
$ws.Range("M105").Value = -8929513
$ws.Range("N105").Value = -5502.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 74917.08
$ws.Range("I31").Value = 5218.5713
$ws.Range("J31").Value = 156232
$ws.Range("K31").Value = 5218.5713
$ws.Range("L31").Value = 156232
$ws.Range("M31").Value = -4923.5713
$ws.Range("N31").Value = -156822
$ws.Range("H34").Value = 74917.08
$ws.Range("I34").Value = 5218.5713
$ws.Range("J34").Value = 156232
$ws.Range("K34").Value = 5218.5713
$ws.Range("L34").Value = 156232
$ws.Range("M34").Value = -5016.5713
$ws.Range("N34").Value = -156636
$ws.Range("H134").Value = 2608.125
$ws.Range("I134").Value = 1791
$ws.Range("K134").Value = 5373
$ws.Range("M134").Value = -2838
$ws.Range("H141").Value = 208064.84
$ws.Range("J141").Value = 208064.84
$ws.Range("L141").Value = 208064.84
$ws.Range("N141").Value = -218424.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 31423.828
$ws.Range("I12").Value = 55601.875
$ws.Range("K12").Value = 166805.625
$ws.Range("M12").Value = -166632.625
$ws.Range("H55").Value = 90911576
$ws.Range("I55").Value = 499500200
$ws.Range("J55").Value = 114100
$ws.Range("K55").Value = 1498500600
$ws.Range("L55").Value = 342300
$ws.Range("M55").Value = -1498500423
$ws.Range("N55").Value = -342654
$ws.Range("H56").Value = 22733672
$ws.Range("I56").Value = 22733672
$ws.Range("K56").Value = 22733672
$ws.Range("M56").Value = -22733142
$ws.Range("H130").Value = 2455.3333
$ws.Range("I130").Value = 2276.75
$ws.Range("J130").Value = 2812.5
$ws.Range("K130").Value = 6830.25
$ws.Range("L130").Value = 8437.5
$ws.Range("M130").Value = -1810.25
$ws.Range("N130").Value = -18477.5
$ws.Range("H131").Value = 10171096
$ws.Range("I131").Value = 7578140.5
$ws.Range("J131").Value = 11121846
$ws.Range("K131").Value = 22734421.5
$ws.Range("L131").Value = 33365538
$ws.Range("M131").Value = -22729381.5
$ws.Range("N131").Value = -33375618
$ws.Range("H137").Value = 5999.5
$ws.Range("I137").Value = 5999
$ws.Range("K137").Value = 17997
$ws.Range("M137").Value = -12897

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 196192.73
$ws.Range("I122").Value = 249471.72
$ws.Range("K122").Value = 748415.16
$ws.Range("M122").Value = -745965.16

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 3733.7273
$ws.Range("I35").Value = 3164.5
$ws.Range("J35").Value = 4416.8
$ws.Range("K35").Value = 3164.5
$ws.Range("L35").Value = 4416.8
$ws.Range("M35").Value = -2828.5
$ws.Range("N35").Value = -5088.8
$ws.Range("H40").Value = 8604
$ws.Range("I40").Value = 4397.4
$ws.Range("J40").Value = 12810.6
$ws.Range("K40").Value = 4397.4
$ws.Range("L40").Value = 12810.6
$ws.Range("M40").Value = -4261.4
$ws.Range("N40").Value = -13082.6
$ws.Range("H46").Value = 5526.364
$ws.Range("I46").Value = 1499
$ws.Range("J46").Value = 5929.1
$ws.Range("K46").Value = 1499
$ws.Range("L46").Value = 5929.1
$ws.Range("N46").Value = -6305.1
$ws.Range("M46").Value = -1311
$ws.Range("H122").Value = 7725.467
$ws.Range("I122").Value = 6644.5454
$ws.Range("J122").Value = 10698
$ws.Range("K122").Value = 19933.6362
$ws.Range("L122").Value = 32094
$ws.Range("M122").Value = -17483.6362
$ws.Range("N122").Value = -36994
$ws.Range("H136").Value = 147487.5
$ws.Range("I136").Value = 170777.08
$ws.Range("K136").Value = 512331.24
$ws.Range("M136").Value = -509781.24

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9042.733
$ws.Range("J62").Value = 9199.393
$ws.Range("L62").Value = 9199.393
$ws.Range("N62").Value = -10447.393
$ws.Range("H65").Value = 9042.733
$ws.Range("J65").Value = 9199.393
$ws.Range("L65").Value = 45996.965
$ws.Range("N65").Value = -52236.965
$ws.Range("H107").Value = 31251096
$ws.Range("J107").Value = 3077.5715
$ws.Range("L107").Value = 9232.7145
$ws.Range("N107").Value = -13072.7145
$ws.Range("H122").Value = 9979.833000000001
$ws.Range("I122").Value = 9999
$ws.Range("J122").Value = 9976
$ws.Range("K122").Value = 29997
$ws.Range("L122").Value = 29928
$ws.Range("M122").Value = -27547
$ws.Range("N122").Value = -34828
$ws.Range("H132").Value = 24073950
$ws.Range("I132").Value = 28576746
$ws.Range("K132").Value = 85730238
$ws.Range("M132").Value = -85727708
$ws.Range("H136").Value = 3079.8518
$ws.Range("I136").Value = 2851.3171
$ws.Range("J136").Value = 3800.6155
$ws.Range("K136").Value = 8553.951300000001
$ws.Range("L136").Value = 11401.8465
$ws.Range("M136").Value = -6003.951300000001
$ws.Range("N136").Value = -16501.8465
